# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> used by the (only) Slide Master   -- currently "Integral"
#   ppt/theme/theme2.xml  -> used by the Notes Master           -- currently "Office Theme"
#
# The target edit swaps the two themes' content, so that the Slide Master ends
# up using the stock "Office Theme" colour scheme (and the Notes Master ends
# up with the "Integral" colours). The two themes already share an identical
# font scheme (fontScheme) and format scheme (fmtScheme) - only the colour
# scheme (clrScheme) actually differs between them - so re-pointing every
# theme colour slot to the "Office Theme" RGB values reproduces the
# meaningful content of the target theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# ThemeColorSchemeIndex order (1-based): Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink.
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # PowerPoint's COM RGB long is B*65536 + G*256 + R (standard OLE colour order).
    $oleColor = ($b * 65536) + ($g * 256) + $r
    $tcs.Colors($i).RGB = $oleColor
}
